# Updated Rise and Settling Time Measurements
# Adds a new "Task4_StartupBehaviour" worksheet (with rise/settling-time
# data derived from the startup transient capture) after "Task4_Efficiency".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new worksheet as the LAST sheet in the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Task4_StartupBehaviour"

# ---------------------------------------------------------------------
# 2. Column widths (approximate - cosmetic only).
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.7
$ws.Columns.Item(4).ColumnWidth = 12.55
$ws.Columns.Item(5).ColumnWidth = 15.3
$ws.Columns.Item(6).ColumnWidth = 12.55
$ws.Columns.Item(7).ColumnWidth = 21.1
$ws.Columns.Item(8).ColumnWidth = 18.7

# ---------------------------------------------------------------------
# 3. Row 1 - merged-looking header over E:G plus H1.
# ---------------------------------------------------------------------
$ws.Cells.Item(1,5).Value = "Rise Time from 0.1 `$V_{out}`$ to 0.9`$V_{out}`$"
$ws.Cells.Item(1,6).Value = ""
$ws.Cells.Item(1,7).Value = ""
$ws.Cells.Item(1,8).Value = '$\Eta < 0.05 V_{out}$'
$ws.Range("E1:G1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 14.4

# ---------------------------------------------------------------------
# 4. Row 2 - column headers.
# ---------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = '$R_{load}$ ($\Omega$)'
$ws.Cells.Item(2,2).Value = '$V_{out}$ (V)'
$ws.Cells.Item(2,3).Value = '$V_{out, max}$ (V)'
$ws.Cells.Item(2,4).Value = "Overshoot (%)"
$ws.Cells.Item(2,5).Value = "Rise Time (s)"
$ws.Cells.Item(2,6).Value = "FROM"
$ws.Cells.Item(2,7).Value = "TO"
$ws.Cells.Item(2,8).Value = "Settling Time (s)"

# ---------------------------------------------------------------------
# 5. Data rows 3-11 : A (Rload), B (Vout), C (Vout,max).
# ---------------------------------------------------------------------
$A = @{3=2.4700000000000002; 4=12.35; 5=24.7; 6=123.5; 7=148.19999999999999; 8=172.9; 9=197.6; 10=222.3; 11=247}
$B = @{3=88.345399999999998; 4=179.64099999999999; 5=259.55599999999998; 6=377.67; 7=388.214; 8=393.375; 9=397.339; 10=400.47899999999998; 11=403.02699999999999}
$C = @{3=158.83600000000001; 4=233.45500000000001; 5=297.37; 6=387.16800000000001; 7=397.96800000000002; 8=402.20699999999999; 9=405.423; 10=407.952; 11=409.983}

foreach ($r in 3..11) {
    $ws.Cells.Item($r,1).Value = $A[$r]
    $ws.Cells.Item($r,2).Value = $B[$r]
    $ws.Cells.Item($r,3).Value = $C[$r]
}

# Column D - overshoot formula = (C-B)/B ; D3 standalone, D4:D11 one shared group.
$ws.Cells.Item(3,4).Formula = "=(C3-B3)/B3"
$ws.Range("D4:D11").Formula = "=(C4-B4)/B4"

# ---------------------------------------------------------------------
# 6. Column E - rise-time text markers for rows 3-5, numeric for 6-11.
#    (scientific notation literals aren't accepted by the script parser,
#    so every value below is written out in plain decimal form.)
# ---------------------------------------------------------------------
$ws.Cells.Item(3,5).Value = "~~-0.00133999~~"
$ws.Cells.Item(4,5).Value = "~~0.0046~~"
$ws.Cells.Item(5,5).Value = "~~0.00604003~~"
$ws.Cells.Item(6,5).Value = 0.0026602
$ws.Cells.Item(7,5).Value = 0.00228069
$ws.Cells.Item(8,5).Value = 0.0019201
$ws.Cells.Item(9,5).Value = 0.0016203
$ws.Cells.Item(10,5).Value = 0.00139993
$ws.Cells.Item(11,5).Value = 0.00123019

# Column F - FROM values (rows 3-11, all numeric).
$ws.Cells.Item(3,6).Value = 0.00409017
$ws.Cells.Item(4,6).Value = 0.00341025
$ws.Cells.Item(5,6).Value = 0.00211016
$ws.Cells.Item(6,6).Value = 0.000649987
$ws.Cells.Item(7,6).Value = 0.000599527
$ws.Cells.Item(8,6).Value = 0.000549968
$ws.Cells.Item(9,6).Value = 0.000519787
$ws.Cells.Item(10,6).Value = 0.000490101
$ws.Cells.Item(11,6).Value = 0.000469821

# Column G - TO values (rows 3-11, all numeric).
$ws.Cells.Item(3,7).Value = 0.00275017
$ws.Cells.Item(4,7).Value = 0.00801025
$ws.Cells.Item(5,7).Value = 0.00815018
$ws.Cells.Item(6,7).Value = 0.00331019
$ws.Cells.Item(7,7).Value = 0.00288022
$ws.Cells.Item(8,7).Value = 0.00247007
$ws.Cells.Item(9,7).Value = 0.00214008
$ws.Cells.Item(10,7).Value = 0.00189003
$ws.Cells.Item(11,7).Value = 0.00170001

# Column H - settling-time text markers for rows 3-5, numeric for 6-11.
$ws.Cells.Item(3,8).Value = "~~0.0199981~~"
$ws.Cells.Item(4,8).Value = "~~0.0199986~~"
$ws.Cells.Item(5,8).Value = "~~0.0199984~~"
$ws.Cells.Item(6,8).Value = 0.00514785
$ws.Cells.Item(7,8).Value = 0.00414786
$ws.Cells.Item(8,8).Value = 0.0034279
$ws.Cells.Item(9,8).Value = 0.00281791
$ws.Cells.Item(10,8).Value = 0.00234793
$ws.Cells.Item(11,8).Value = 0.00200796

# ---------------------------------------------------------------------
# 7. Selection / activation - new sheet becomes the active tab, with
#    the cursor parked at J24 (matches the source workbook's saved view).
# ---------------------------------------------------------------------
$ws.Range("J24").Select()
$ws.Activate()

$wb.Application.Calculate()
